# PITCH.pptx - "Correção do PITCH versão final - ajuste do LOGO"
#
# The "Logo do Projeto" placeholder (shape 3 on slide 1, originally an
# Oval named "Oval 3") is converted into a Rectangle: renamed, moved,
# resized, its geometry switched from ellipse to rect, and its caption
# text gets a trailing space while the stray empty trailing paragraph
# is dropped.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
try {
    $sh = $s.Shapes.Item("Oval 3")
} catch {
    $sh = $s.Shapes.Item(3)
}

# Rename "Oval 3" -> "Rectangle 3"
$sh.Name = "Rectangle 3"

# Change the autoshape geometry from ellipse to rectangle
$sh.AutoShapeType = 1   # msoShapeRectangle

# Move + resize (points; chosen so the EMU round-trip lands exactly on
# the target off/ext values of x=643467 y=3928533 cx=2551289 cy=1655762)
$sh.Left   = 50.66669464111328
$sh.Top    = 309.33331298828125
$sh.Width  = 200.88890075683594
$sh.Height = 130.37496948242188

# Update the caption text (trailing space added) - this also collapses
# the now-redundant empty trailing paragraph into the single remaining one
$sh.TextFrame.TextRange.Text = "Logo do Projeto "
